# Updated cryptos list on Sat Feb 25 20:01:08 UTC 2023 with GitHub Actions
#
# Refreshes the "Coin / Link / Price / Volume(1h)" table on Sheet1 with the
# latest scrape: new Price/Volume figures for every coin, plus a rank swap
# between Dogecoin and Polygon (rows 11/12 traded places, each taking the
# other's name/link along with its own freshly scraped price + volume).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row's updates: Cell => new text value, exactly as scraped (kept as
# literal text - these "prices" use dotted thousands separators, e.g.
# "23.045.60", and the volumes keep their padded "  +0.33%  " layout).
$updates = @(
    @{ Cell = "D2"; Value = "23.045.60" }
    @{ Cell = "E2"; Value = "  -0.52%  " }
    @{ Cell = "D3"; Value = "1.587.84" }
    @{ Cell = "E3"; Value = "  -1.32%  " }
    @{ Cell = "D4"; Value = "1.002" }
    @{ Cell = "E4"; Value = "  +0.33%  " }
    @{ Cell = "D5"; Value = "1.002" }
    @{ Cell = "E5"; Value = "  +0.27%  " }
    @{ Cell = "D6"; Value = "301.03" }
    @{ Cell = "E6"; Value = "  -0.34%  " }
    @{ Cell = "D7"; Value = "0.3768" }
    @{ Cell = "E7"; Value = "  -0.19%  " }
    @{ Cell = "D8"; Value = "0.3591" }
    @{ Cell = "E8"; Value = "  -2.00%  " }
    @{ Cell = "D9"; Value = "50.46" }
    @{ Cell = "E9"; Value = "  +3.47%  " }
    @{ Cell = "D10"; Value = "1.003" }
    @{ Cell = "E10"; Value = "  +0.37%  " }
    @{ Cell = "B11"; Value = "Polygon" }
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic" }
    @{ Cell = "D11"; Value = "1.219" }
    @{ Cell = "E11"; Value = "  -4.64%  " }
    @{ Cell = "B12"; Value = "Dogecoin" }
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge" }
    @{ Cell = "D12"; Value = "0.08041" }
    @{ Cell = "E12"; Value = "  -0.67%  " }
    @{ Cell = "D13"; Value = "21.98" }
    @{ Cell = "E13"; Value = "  -4.58%  " }
    @{ Cell = "D14"; Value = "6.485" }
    @{ Cell = "E14"; Value = "  -2.36%  " }
    @{ Cell = "D15"; Value = "7.338" }
    @{ Cell = "E15"; Value = "  -4.38%  " }
    @{ Cell = "D16"; Value = "0.00001224" }
    @{ Cell = "E16"; Value = "  -3.58%  " }
    @{ Cell = "D17"; Value = "1.588.80" }
    @{ Cell = "E17"; Value = "  -0.57%  " }
    @{ Cell = "D18"; Value = "92.22" }
    @{ Cell = "D19"; Value = "0.06786" }
    @{ Cell = "E19"; Value = "  -0.13%  " }
    @{ Cell = "D20"; Value = "17.92" }
    @{ Cell = "E20"; Value = "  -2.79%  " }
    @{ Cell = "E21"; Value = "  +0.17%  " }
    @{ Cell = "D22"; Value = "6.432" }
    @{ Cell = "E22"; Value = "  -2.45%  " }
    @{ Cell = "D23"; Value = "12.81" }
    @{ Cell = "E23"; Value = "  -1.73%  " }
    @{ Cell = "D24"; Value = "23.010.72" }
    @{ Cell = "E24"; Value = "  -0.65%  " }
    @{ Cell = "E25"; Value = "  +0.63%  " }
    @{ Cell = "D26"; Value = "2.763" }
    @{ Cell = "E26"; Value = "  -5.99%  " }
    @{ Cell = "D27"; Value = "20.81" }
    @{ Cell = "E27"; Value = "  -1.49%  " }
    @{ Cell = "D28"; Value = "147.40" }
    @{ Cell = "E28"; Value = "  -2.24%  " }
    @{ Cell = "D29"; Value = "5.203" }
    @{ Cell = "E29"; Value = "  -1.36%  " }
    @{ Cell = "D30"; Value = "133.00" }
    @{ Cell = "E30"; Value = "  +0.53%  " }
    @{ Cell = "D31"; Value = "2.319" }
    @{ Cell = "E31"; Value = "  -3.68%  " }
    @{ Cell = "D32"; Value = "6.537" }
    @{ Cell = "E32"; Value = "  -5.96%  " }
    @{ Cell = "D33"; Value = "1.763.94" }
    @{ Cell = "E33"; Value = "  -0.52%  " }
    @{ Cell = "D34"; Value = "0.9413" }
    @{ Cell = "E34"; Value = "  -5.36%  " }
    @{ Cell = "D35"; Value = "0.07340" }
    @{ Cell = "E35"; Value = "  -5.14%  " }
    @{ Cell = "D36"; Value = "0.02678" }
    @{ Cell = "E36"; Value = "  -4.27%  " }
    @{ Cell = "D37"; Value = "10.04" }
    @{ Cell = "E37"; Value = "  -1.22%  " }
    @{ Cell = "D38"; Value = "0.08763" }
    @{ Cell = "E38"; Value = "  -1.17%  " }
    @{ Cell = "D39"; Value = "6.079" }
    @{ Cell = "E39"; Value = "  -3.67%  " }
    @{ Cell = "D40"; Value = "0.2483" }
    @{ Cell = "E40"; Value = "  -2.61%  " }
    @{ Cell = "D41"; Value = "1.337" }
    @{ Cell = "E41"; Value = "  -4.50%  " }
    @{ Cell = "D42"; Value = "0.6884" }
    @{ Cell = "E42"; Value = "  -4.11%  " }
    @{ Cell = "D43"; Value = "11.92" }
    @{ Cell = "E43"; Value = "  -7.01%  " }
    @{ Cell = "D44"; Value = "14.86" }
    @{ Cell = "E44"; Value = "  -6.91%  " }
    @{ Cell = "D45"; Value = "0.6409" }
    @{ Cell = "E45"; Value = "  -3.35%  " }
    @{ Cell = "D46"; Value = "3.989" }
    @{ Cell = "E46"; Value = "  +0.50%  " }
    @{ Cell = "D47"; Value = "2.248" }
    @{ Cell = "E47"; Value = "  -2.90%  " }
    @{ Cell = "D48"; Value = "131.51" }
    @{ Cell = "E48"; Value = "  -0.04%  " }
    @{ Cell = "D49"; Value = "0.07883" }
    @{ Cell = "E49"; Value = "  -1.80%  " }
    @{ Cell = "D50"; Value = "1.197" }
    @{ Cell = "E50"; Value = "  +1.96%  " }
    @{ Cell = "D51"; Value = "1.189" }
    @{ Cell = "E51"; Value = "  +0.55%  " }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $val = $u.Value

    # The sheet stores these figures as literal text (inline strings), not
    # numbers - e.g. "23.045.60" uses a dotted thousands separator and isn't
    # a real number at all, while values like "301.03" or "133.00" *do*
    # look numeric to Excel. Left alone, Range.Value would silently convert
    # those into floating-point numbers (and mangle trailing zeros / introduce
    # rounding noise), so force text formatting first whenever the new value
    # parses as a plain number, keeping the digits exactly as scraped.
    if ($val -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $rng.NumberFormat = "@"
    }

    $rng.Value = $val
}
